$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.387.78"
$ws.Range("E2").Value = "  +4.20%  "
$ws.Range("D3").Value = "2.471.51"
$ws.Range("E3").Value = "  +2.43%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'569.38"
$ws.Range("E5").Value = "  +2.74%  "
$ws.Range("D6").Value = "'169.01"
$ws.Range("E6").Value = "  +5.35%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("E9").Value = "  +14.36%  "
$ws.Range("D10").Value = "2.468.45"
$ws.Range("E10").Value = "  +2.24%  "
$ws.Range("E11").Value = "  -1.38%  "
$ws.Range("D12").Value = "'0.338"
$ws.Range("E12").Value = "  +4.02%  "
$ws.Range("D13").Value = "'4.73"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("E14").Value = "  +9.91%  "
$ws.Range("D15").Value = "70.250.28"
$ws.Range("E15").Value = "  +4.05%  "
$ws.Range("D16").Value = "2.920.20"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").Value = "'24.36"
$ws.Range("E17").Value = "  +6.36%  "
$ws.Range("D18").Value = "2.459.92"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "'10.92"
$ws.Range("E19").Value = "  +6.10%  "
$ws.Range("D20").Value = "'7.23"
$ws.Range("E20").Value = "  +6.23%  "
$ws.Range("D21").Value = "'344.88"
$ws.Range("E21").Value = "  +3.29%  "
$ws.Range("D22").Value = "'3.91"
$ws.Range("E22").Value = "  +3.94%  "
$ws.Range("D23").Value = "'2.04"
$ws.Range("E23").Value = "  +9.58%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "'66.79"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").Value = "'3.92"
$ws.Range("E26").Value = "  +9.05%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "'1.10"
$ws.Range("E27").Value = "  +10.35%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "'8.62"
$ws.Range("E28").Value = "  +7.66%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.586.53"
$ws.Range("E29").Value = "  +1.76%  "
$ws.Range("D30").Value = "0.0₃0873"
$ws.Range("E30").Value = "  +8.72%  "
$ws.Range("D31").Value = "'7.43"
$ws.Range("E31").Value = "  +5.08%  "
$ws.Range("E32").Value = "  +12.36%  "
$ws.Range("D33").Value = "'459.68"
$ws.Range("E33").Value = "  +9.64%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("E35").Value = "  +3.01%  "
$ws.Range("D36").Value = "'161.84"
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.113"
$ws.Range("E37").Value = "  +10.89%  "
$ws.Range("B38").Value = "WhiteBITCoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D38").Value = "'19.12"
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("D40").Value = "'18.32"
$ws.Range("E40").Value = "  +3.71%  "
$ws.Range("D41").Value = "'0.307"
$ws.Range("E41").Value = "  +5.19%  "
$ws.Range("D42").Value = "'1.56"
$ws.Range("E42").Value = "  +7.52%  "
$ws.Range("D43").Value = "'4.47"
$ws.Range("E43").Value = "  +5.42%  "
$ws.Range("D44").Value = "'37.92"
$ws.Range("E44").Value = "  +1.22%  "
$ws.Range("D45").Value = "'1.11"
$ws.Range("E45").Value = "  +6.04%  "
$ws.Range("D46").Value = "'2.17"
$ws.Range("E46").Value = "  +9.02%  "
$ws.Range("D47").Value = "'3.43"
$ws.Range("E47").Value = "  +3.22%  "
$ws.Range("D48").Value = "'134.47"
$ws.Range("E48").Value = "  +5.12%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "'0.496"
$ws.Range("E49").Value = "  +4.64%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0727"
$ws.Range("E50").Value = "  +2.44%  "
$ws.Range("D51").Value = "'0.567"
$ws.Range("E51").Value = "  +2.45%  "
